# Fix the product name string (insert missing dash) on both sheets, and
# switch the active sheet/selection from ProductLoanInput to ProductLoanOutput.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newName = "343-MS-EPP-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"

# Correct the product name text on both sheets (missing dash after "343").
$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# Update the remembered selection on the input sheet (it is no longer the
# active tab, but Excel still persists its last selection).
$wsInput.Activate()
$wsInput.Range("B1").Select()

# Make the output sheet the active/selected tab with B1 selected.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
